$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G (shifts existing G:N right to H:O).
# CopyOrigin 0 = xlFormatFromLeftOrAbove, so the new column inherits
# formatting from column F (matches Excel's default Insert behavior).
$ws.Columns("G:G").Insert($null, 0)

# The new column keeps the same width as its left neighbour (F).
$ws.Columns("G:G").ColumnWidth = $ws.Columns("F:F").ColumnWidth

# New "GASTO" header cell, formatted like the other header cells (e.g. B9).
$ws.Range("B9").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("G9").Value = "GASTO"

$excel.CutCopyMode = 0
$ws.Range("G10").Select() | Out-Null
